$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

# Columns whose text content looks like a pure number/date would otherwise be
# auto-converted by Excel's type inference. Prefix with an apostrophe to force
# text entry (this is how a human typing into Excel would keep e.g. "6076" or
# "6/24/2025" as text), then reset the cell style to "Normal" so the
# quote-prefix formatting doesn't leave a visible style index on the cell.
function Set-TextValue($cell, [string]$text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item($row, 1) "6076"
Set-TextValue $ws.Cells.Item($row, 2) "6/24/2025"
$ws.Cells.Item($row, 3).Value = "MATHEU 727"
Set-TextValue $ws.Cells.Item($row, 4) "3"
Set-TextValue $ws.Cells.Item($row, 5) ""
$ws.Cells.Item($row, 6).Value = "GESTION TELECENTRO"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Picada"
Set-TextValue $ws.Cells.Item($row, 9) "1"
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Nodo TLC"
$ws.Cells.Item($row, 12).Value = "Pasante"
$ws.Cells.Item($row, 13).Value = -58.400169
$ws.Cells.Item($row, 14).Value = -34.617784
